$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# K2 / K3 swap
$ws.Range("K2").Value = 0
$ws.Range("K3").Value = 1

# Row 6 (was Biagio VDS) -> now Roy Droog's stats
$ws.Range("B6").Value = "Roy Droog"
$ws.Range("C6").Value = 12
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6984
$ws.Range("H6").Value = 410
$ws.Range("I6").Value = 51.1

# Row 7 (was Ewan Taylor) -> now Biagio VDS's stats
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Biagio VDS"
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 6181
$ws.Range("H7").Value = 269
$ws.Range("I7").Value = 68.93000000000001
$ws.Range("J7").Value = 13

# Row 8 (was wessel de Haas) -> now Ewan Taylor's stats
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Ewan Taylor"
$ws.Range("G8").Value = 6391
$ws.Range("H8").Value = 337
$ws.Range("I8").Value = 56.89

# Row 9 (was Robin Willis) -> now wessel de Haas's stats
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "wessel de Haas"
$ws.Range("C9").Value = 11
$ws.Range("G9").Value = 9005
$ws.Range("H9").Value = 558
$ws.Range("I9").Value = 48.41
$ws.Range("J9").Value = 11

# Row 10 (was Roy Droog) -> now Robin Willis's stats
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Robin Willis"
$ws.Range("C10").Value = 10
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 4671
$ws.Range("H10").Value = 279
$ws.Range("I10").Value = 50.23
